# Rename the sole worksheet from "Sheet1" to "Mitra Oddity".
# (Excel automatically re-points any defined names / formulas that
# reference the sheet, e.g. the _xlnm._FilterDatabase name below.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Mitra Oddity"

# The workbook accumulated a large number of stale/legacy defined names
# (old #REF!, NA(), #N/A leftovers from long-dead external links). Strip
# all of them, keeping only the built-in _FilterDatabase name that the
# AutoFilter on the sheet still needs.
$guard = 0
while ($wb.Names.Count -gt 1 -and $guard -lt 1000) {
    $guard += 1
    $deletedOne = $false
    for ($i = 1; $i -le $wb.Names.Count; $i++) {
        $n = $wb.Names.Item($i)
        if ($n.Name -notlike "*_FilterDatabase*") {
            $n.Delete()
            $deletedOne = $true
            break
        }
    }
    if (-not $deletedOne) {
        break
    }
}
